$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column C slightly to fit the new longer example text
# (closest achievable raw width to the authored 33.6640625 given
# the engine's internal pixel-width quantization)
$ws.Columns.Item(3).ColumnWidth = 32.83

# Add new example rows (3 and 4) with formatter strings
$ws.Range("B3").Value = "number(data.numberVal 0%)"
$ws.Range("C3").Value = "date(data.dateVal dd.mm.yyyy hh:mm)"
$ws.Range("C4").Value = "date(data.dateVal dd.mm.yyyy)"

# Update the active selection to match the authored state
$ws.Range("C6").Select()
